# Weekly update: insert a new record (row) for the Cilantro market data,
# pushing existing rows 6-15 down to 7-16, and populate the new row 6
# with the latest observation (Fecha=44690, Volumen=500) while the rest
# of the fields repeat the ones from the previous top entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 6; this shifts old rows 6..15 to 7..16
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with the new data point
$ws.Cells.Item(6, 1).Value  = 5
$ws.Cells.Item(6, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(6, 3).Value  = "Maule"
$ws.Cells.Item(6, 4).Value  = 44690
$ws.Cells.Item(6, 5).Value  = 7
$ws.Cells.Item(6, 6).Value  = 100112040
$ws.Cells.Item(6, 7).Value  = "Cilantro"
$ws.Cells.Item(6, 8).Value  = "Sin especificar"
$ws.Cells.Item(6, 9).Value  = "Primera"
$ws.Cells.Item(6, 10).Value = 500
$ws.Cells.Item(6, 11).Value = 7000
$ws.Cells.Item(6, 12).Value = 7000
$ws.Cells.Item(6, 13).Value = 7000
$ws.Cells.Item(6, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(6, 15).Value = "Región del Maule"
$ws.Cells.Item(6, 16).Value = 194
$ws.Cells.Item(6, 17).Value = 36
$ws.Cells.Item(6, 18).Value = "Hortaliza"
